# Apply the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Tue Aug 29 04:44:04 UTC 2023 with GitHub Actions".
# Column D = Price (text, e.g. "218.78"), Column E = Volume(1h) (text, e.g. "  +0.05%  ").
#
# Many Price strings look like plain numbers ("218.78", "20.58", ...). Assigning
# them straight to .Value would make Excel auto-coerce the cell to a Number
# (losing the original text formatting/precision, e.g. "0.5238" -> 0.52380000000000004).
# To keep them as text - matching the workbook's original inlineStr cells - we force
# the cell to Text format before writing, then reset the style so no stray
# per-cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '26.151.88'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.656.08'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.19%  '
Set-TextValue "D5" '218.78'
$ws.Range("E5").Value = '  +0.05%  '
Set-TextValue "D6" '0.5238'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E7").Value = '  -0.16%  '
Set-TextValue "D8" '0.2658'
$ws.Range("E8").Value = '  +1.41%  '
Set-TextValue "D9" '0.06357'
$ws.Range("E9").Value = '  +1.03%  '
Set-TextValue "D10" '20.58'
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("E11").Value = '  -1.47%  '
Set-TextValue "D12" '4.604'
$ws.Range("E12").Value = '  +2.37%  '
$ws.Range("D13").Value = '1.645.91'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '1.884.98'
$ws.Range("E14").Value = '  +0.05%  '
Set-TextValue "D15" '0.5625'
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").Value = '0.0₅8197'
$ws.Range("E16").Value = '  +2.45%  '
Set-TextValue "D17" '65.47'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '26.146.66'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("E19").Value = '  -0.17%  '
Set-TextValue "D20" '4.659'
$ws.Range("E20").Value = '  +0.46%  '
Set-TextValue "D21" '10.52'
$ws.Range("E21").Value = '  +4.03%  '
Set-TextValue "D22" '192.38'
$ws.Range("E22").Value = '  -1.58%  '
Set-TextValue "D23" '5.961'
$ws.Range("E23").Value = '  -0.02%  '
Set-TextValue "D24" '1.003'
$ws.Range("E24").Value = '  -0.25%  '
Set-TextValue "D25" '145.23'
$ws.Range("E25").Value = '  -0.98%  '
Set-TextValue "D26" '0.1197'
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("E27").Value = '  +1.45%  '
Set-TextValue "D28" '15.96'
$ws.Range("E28").Value = '  +0.07%  '
Set-TextValue "D29" '1.513'
$ws.Range("E29").Value = '  +1.00%  '
Set-TextValue "D30" '0.05467'
$ws.Range("E30").Value = '  -4.42%  '
Set-TextValue "D31" '1.274'
$ws.Range("E31").Value = '  +0.21%  '
Set-TextValue "D32" '3.467'
$ws.Range("E32").Value = '  -0.64%  '
Set-TextValue "D33" '3.371'
$ws.Range("E33").Value = '  +0.58%  '
Set-TextValue "D34" '1.565'
$ws.Range("E34").Value = '  -1.39%  '
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("E36").Value = '  -0.87%  '
Set-TextValue "D37" '2.401'
$ws.Range("E37").Value = '  -0.66%  '
Set-TextValue "D38" '0.5683'
$ws.Range("E38").Value = '  -0.17%  '
Set-TextValue "D39" '0.01587'
$ws.Range("E39").Value = '  -0.59%  '
Set-TextValue "D40" '5.881'
$ws.Range("E40").Value = '  -1.23%  '
$ws.Range("E41").Value = '  -0.18%  '
Set-TextValue "D42" '0.8323'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").Value = '1.028.61'
$ws.Range("E43").Value = '  -2.84%  '
Set-TextValue "D44" '101.20'
$ws.Range("E44").Value = '  -2.19%  '
$ws.Range("D45").Value = '1.795.73'
$ws.Range("E45").Value = '  +0.06%  '
Set-TextValue "D46" '57.79'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").Value = '0.0₈107'
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("E48").Value = '  -0.64%  '
Set-TextValue "D49" '8.040'
$ws.Range("E49").Value = '  +0.25%  '
Set-TextValue "D50" '0.4344'
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("E51").Value = '  -3.83%  '
